$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one
# day (2023-09-19 -> 2023-09-20, serial 45188 -> 45189) for every data row
# (rows 2 through 349).
$ws.Range("C2:C349").Value = 45189
